# Auto-applies the numeric-value corrections from the commit diff
# (price/profit recalculations across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets).
$wb = $excel.ActiveWorkbook

# ALC!row 11
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 769285.7
$ws.Range("I11").Value = 769285.7
$ws.Range("K11").Value = 769285.7
$ws.Range("M11").Value = -769145.7

# ALC!row 33
$ws.Range("H33").Value = 33374076
$ws.Range("I33").Value = 40000860
$ws.Range("J33").Value = 240159.2
$ws.Range("K33").Value = 40000860
$ws.Range("L33").Value = 240159.2
$ws.Range("M33").Value = -40000631
$ws.Range("N33").Value = -240617.2

# ALC!row 40
$ws.Range("H40").Value = 4724.7354
$ws.Range("I40").Value = 7219.9414
$ws.Range("J40").Value = 2229.5293
$ws.Range("K40").Value = 7219.9414
$ws.Range("L40").Value = 2229.5293
$ws.Range("M40").Value = -7044.9414
$ws.Range("N40").Value = -2579.5293

# ALC!row 99
$ws.Range("H99").Value = 1045.5555
$ws.Range("I99").Value = 866.2
$ws.Range("J99").Value = 1269.75
$ws.Range("K99").Value = 2598.6
$ws.Range("L99").Value = 3809.25
$ws.Range("M99").Value = -1100.6
$ws.Range("N99").Value = -6805.25

# ALC!row 116
$ws.Range("H116").Value = 5430
$ws.Range("I116").Value = 2265.2942
$ws.Range("J116").Value = 8120
$ws.Range("K116").Value = 2265.2942
$ws.Range("L116").Value = 8120
$ws.Range("M116").Value = 1176.7058
$ws.Range("N116").Value = -15004

# ALC!row 137
$ws.Range("H137").Value = 846.7895
$ws.Range("I137").Value = 739.2
$ws.Range("K137").Value = 2217.6
$ws.Range("M137").Value = 332.3999999999996

# ARM!row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2068.1052
$ws.Range("I61").Value = 1819.6
$ws.Range("K61").Value = 1819.6
$ws.Range("M61").Value = -1607.6

# ARM!row 132
$ws.Range("H132").Value = 1911.0588
$ws.Range("I132").Value = 1436.6154
$ws.Range("J132").Value = 3453
$ws.Range("K132").Value = 4309.8462
$ws.Range("L132").Value = 10359
$ws.Range("M132").Value = -1779.8462
$ws.Range("N132").Value = -15419

# ARM!row 136
$ws.Range("H136").Value = 2068.1052
$ws.Range("I136").Value = 1819.6
$ws.Range("K136").Value = 5458.799999999999
$ws.Range("M136").Value = -2908.799999999999

# BSM!row 10
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 105
$ws.Range("I10").Value = 105
$ws.Range("K10").Value = 105
$ws.Range("M10").Value = 35

# BSM!row 24
$ws.Range("H24").Value = 20050
$ws.Range("I24").Value = 650
$ws.Range("J24").Value = 29750
$ws.Range("K24").Value = 650
$ws.Range("L24").Value = 29750
$ws.Range("M24").Value = -415
$ws.Range("N24").Value = -30220

# BSM!row 86
$ws.Range("H86").Value = 2032.1538
$ws.Range("I86").Value = 1603
$ws.Range("J86").Value = 2222.889
$ws.Range("K86").Value = 1603
$ws.Range("L86").Value = 2222.889
$ws.Range("M86").Value = -480
$ws.Range("N86").Value = -4468.889

# BSM!row 89
$ws.Range("H89").Value = 2032.1538
$ws.Range("I89").Value = 1603
$ws.Range("J89").Value = 2222.889
$ws.Range("K89").Value = 8015
$ws.Range("L89").Value = 11114.445
$ws.Range("M89").Value = -2399
$ws.Range("N89").Value = -22346.445

# BSM!row 107
$ws.Range("H107").Value = 1230.138
$ws.Range("I107").Value = 1245.1538
$ws.Range("J107").Value = 1100
$ws.Range("K107").Value = 1245.1538
$ws.Range("L107").Value = 1100
$ws.Range("M107").Value = 674.8462
$ws.Range("N107").Value = -4940

# CRP!row 2
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 27501.25
$ws.Range("I2").Value = 13333.333
$ws.Range("K2").Value = 13333.333
$ws.Range("M2").Value = -13220.333

# CRP!row 3
$ws.Range("H3").Value = 170133.33
$ws.Range("I3").Value = 250200
$ws.Range("K3").Value = 250200
$ws.Range("M3").Value = -250087

# CRP!row 4
$ws.Range("H4").Value = 2864571.8
$ws.Range("J4").Value = 3341667
$ws.Range("L4").Value = 3341667
$ws.Range("N4").Value = -3341891

# CRP!row 5
$ws.Range("H5").Value = 1771.1428
$ws.Range("I5").Value = 147.5
$ws.Range("J5").Value = 3936
$ws.Range("K5").Value = 147.5
$ws.Range("L5").Value = 3936
$ws.Range("M5").Value = -35.5
$ws.Range("N5").Value = -4160

# CRP!row 11
$ws.Range("H11").Value = 10786.571
$ws.Range("J11").Value = 10786.571
$ws.Range("L11").Value = 10786.571
$ws.Range("N11").Value = -11066.571

# CRP!row 16
$ws.Range("H16").Value = 4344.48
$ws.Range("I16").Value = 3886.7334
$ws.Range("J16").Value = 5031.1
$ws.Range("K16").Value = 3886.7334
$ws.Range("L16").Value = 5031.1
$ws.Range("M16").Value = -3599.7334
$ws.Range("N16").Value = -5605.1

# CRP!row 31
$ws.Range("H31").Value = 1506.9767
$ws.Range("I31").Value = 1098.5927
$ws.Range("J31").Value = 2196.125
$ws.Range("K31").Value = 1098.5927
$ws.Range("L31").Value = 2196.125
$ws.Range("M31").Value = -803.5926999999999
$ws.Range("N31").Value = -2786.125

# CRP!row 34
$ws.Range("H34").Value = 1506.9767
$ws.Range("I34").Value = 1098.5927
$ws.Range("J34").Value = 2196.125
$ws.Range("K34").Value = 1098.5927
$ws.Range("L34").Value = 2196.125
$ws.Range("M34").Value = -896.5926999999999
$ws.Range("N34").Value = -2600.125

# CRP!row 107
$ws.Range("H107").Value = 1220.6666
$ws.Range("I107").Value = 689.2308
$ws.Range("J107").Value = 2084.25
$ws.Range("K107").Value = 689.2308
$ws.Range("L107").Value = 2084.25
$ws.Range("M107").Value = 1230.7692
$ws.Range("N107").Value = -5924.25

# CRP!row 113
$ws.Range("H113").Value = 4344.48
$ws.Range("I113").Value = 3886.7334
$ws.Range("J113").Value = 5031.1
$ws.Range("K113").Value = 3886.7334
$ws.Range("L113").Value = 5031.1
$ws.Range("M113").Value = -1716.7334
$ws.Range("N113").Value = -9371.1

# CRP!row 134
$ws.Range("H134").Value = 1302.2667
$ws.Range("I134").Value = 1222.9412
$ws.Range("J134").Value = 1406
$ws.Range("K134").Value = 3668.8236
$ws.Range("L134").Value = 4218
$ws.Range("M134").Value = -1133.8236
$ws.Range("N134").Value = -9288

# CUL!row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 19.352942
$ws.Range("J12").Value = 23.90909
$ws.Range("L12").Value = 71.72727
$ws.Range("N12").Value = -417.72727

# CUL!row 38
$ws.Range("H38").Value = 83.42856999999999
$ws.Range("I38").Value = 46.444443
$ws.Range("J38").Value = 150
$ws.Range("K38").Value = 139.333329
$ws.Range("L38").Value = 450
$ws.Range("M38").Value = 207.666671
$ws.Range("N38").Value = -1144

# GSM!row 18
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 11944.25
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 11944.25
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 11944.25
$ws.Range("M18").Value = ""
$ws.Range("N18").Value = -12530.25

# GSM!row 43
$ws.Range("H43").Value = 10202.429
$ws.Range("I43").Value = 2283.4
$ws.Range("J43").Value = 30000
$ws.Range("K43").Value = 2283.4
$ws.Range("L43").Value = 30000
$ws.Range("M43").Value = -2132.4
$ws.Range("N43").Value = -30302

# GSM!row 46
$ws.Range("H46").Value = 35330
$ws.Range("J46").Value = 47995
$ws.Range("L46").Value = 47995
$ws.Range("N46").Value = -48307

# GSM!row 70
$ws.Range("H70").Value = 5927
$ws.Range("I70").Value = 5529.75
$ws.Range("J70").Value = 6293.6924
$ws.Range("K70").Value = 5529.75
$ws.Range("L70").Value = 6293.6924
$ws.Range("M70").Value = -5259.75
$ws.Range("N70").Value = -6833.6924

# GSM!row 73
$ws.Range("H73").Value = 5927
$ws.Range("I73").Value = 5529.75
$ws.Range("J73").Value = 6293.6924
$ws.Range("K73").Value = 5529.75
$ws.Range("L73").Value = 6293.6924
$ws.Range("M73").Value = -4593.75
$ws.Range("N73").Value = -8165.6924

# GSM!row 80
$ws.Range("H80").Value = 2266.6667
$ws.Range("I80").Value = 2300
$ws.Range("J80").Value = 2233.3333
$ws.Range("K80").Value = 2300
$ws.Range("L80").Value = 2233.3333
$ws.Range("M80").Value = -1302
$ws.Range("N80").Value = -4229.3333

# GSM!row 83
$ws.Range("H83").Value = 2266.6667
$ws.Range("I83").Value = 2300
$ws.Range("J83").Value = 2233.3333
$ws.Range("K83").Value = 11500
$ws.Range("L83").Value = 11166.6665
$ws.Range("M83").Value = -6508
$ws.Range("N83").Value = -21150.6665

# LTW!row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 421
$ws.Range("I22").Value = 367.10526
$ws.Range("J22").Value = 494.14285
$ws.Range("K22").Value = 367.10526
$ws.Range("L22").Value = 494.14285
$ws.Range("M22").Value = -72.10525999999999
$ws.Range("N22").Value = -1084.14285

# LTW!row 27
$ws.Range("H27").Value = 421
$ws.Range("I27").Value = 367.10526
$ws.Range("J27").Value = 494.14285
$ws.Range("K27").Value = 367.10526
$ws.Range("L27").Value = 494.14285
$ws.Range("M27").Value = -260.10526
$ws.Range("N27").Value = -708.14285

# LTW!row 46
$ws.Range("H46").Value = 1199.8889
$ws.Range("I46").Value = 1299.75
$ws.Range("J46").Value = 1120
$ws.Range("K46").Value = 1299.75
$ws.Range("L46").Value = 1120
$ws.Range("M46").Value = -1111.75
$ws.Range("N46").Value = -1496

# LTW!row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""

# WVR!row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7588.684
$ws.Range("I81").Value = 14748.875
$ws.Range("J81").Value = 2381.2727
$ws.Range("K81").Value = 29497.75
$ws.Range("L81").Value = 4762.5454
$ws.Range("M81").Value = -28436.75
$ws.Range("N81").Value = -6884.5454

# WVR!row 84
$ws.Range("H84").Value = 7588.684
$ws.Range("I84").Value = 14748.875
$ws.Range("J84").Value = 2381.2727
$ws.Range("K84").Value = 147488.75
$ws.Range("L84").Value = 23812.727
$ws.Range("M84").Value = -142184.75
$ws.Range("N84").Value = -34420.727
